$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9771034717559814
$ws.Range("B1").Value = 1.567568898200989
$ws.Range("D1").Value = 1.767438530921936
$ws.Range("E1").Value = 1.057743430137634
